# Implemented getting kafka relations.
# The reflection-derived "classFields" sheet was regenerated and several
# classes now list their fields in a different order:
#   - com.macro.mall.auth.constant.MessageConstant
#   - com.macro.mall.auth.config.Oauth2ServerConfig
#   - com.macro.mall.auth.domain.Oauth2TokenDto
#   - com.macro.mall.auth.domain.Oauth2TokenDto$Oauth2TokenDtoBuilder
#   - com.macro.mall.auth.domain.SecurityUser
# This script rewrites the affected rows (B = Field Name, D = Field Type)
# on the classFields worksheet to match the new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# com.macro.mall.auth.constant.MessageConstant (rows 2-8): field name only
$ws.Range("B2").Value = "USERNAME_PASSWORD_ERROR"
$ws.Range("B3").Value = "ACCOUNT_DISABLED"
$ws.Range("B4").Value = "ACCOUNT_EXPIRED"
$ws.Range("B5").Value = "LOGIN_SUCCESS"
$ws.Range("B6").Value = "ACCOUNT_LOCKED"
$ws.Range("B7").Value = "PERMISSION_DENIED"
$ws.Range("B8").Value = "CREDENTIALS_EXPIRED"

# com.macro.mall.auth.config.Oauth2ServerConfig (rows 10-13)
$ws.Range("B10").Value = "authenticationManager"
$ws.Range("D10").Value = "org.springframework.security.authentication.AuthenticationManager"
$ws.Range("B11").Value = "passwordEncoder"
$ws.Range("D11").Value = "org.springframework.security.crypto.password.PasswordEncoder"
$ws.Range("B12").Value = "userDetailsService"
$ws.Range("D12").Value = "com.macro.mall.auth.service.impl.UserServiceImpl"
$ws.Range("B13").Value = "jwtTokenEnhancer"
$ws.Range("D13").Value = "com.macro.mall.auth.component.JwtTokenEnhancer"

# com.macro.mall.auth.domain.Oauth2TokenDto (rows 14-17)
$ws.Range("B14").Value = "expiresIn"
$ws.Range("D14").Value = "int"
$ws.Range("B15").Value = "tokenHead"
$ws.Range("D15").Value = "java.lang.String"
$ws.Range("B16").Value = "token"
$ws.Range("D16").Value = "java.lang.String"
$ws.Range("B17").Value = "refreshToken"
$ws.Range("D17").Value = "java.lang.String"

# com.macro.mall.auth.domain.Oauth2TokenDto$Oauth2TokenDtoBuilder (rows 21-24)
$ws.Range("B21").Value = "token"
$ws.Range("D21").Value = "java.lang.String"
$ws.Range("B22").Value = "expiresIn"
$ws.Range("D22").Value = "int"
$ws.Range("B23").Value = "refreshToken"
$ws.Range("D23").Value = "java.lang.String"
$ws.Range("B24").Value = "tokenHead"
$ws.Range("D24").Value = "java.lang.String"

# com.macro.mall.auth.domain.SecurityUser (rows 25-30)
$ws.Range("B25").Value = "password"
$ws.Range("D25").Value = "java.lang.String"
$ws.Range("B26").Value = "authorities"
$ws.Range("D26").Value = "java.util.Collection"
$ws.Range("B27").Value = "clientId"
$ws.Range("D27").Value = "java.lang.String"
$ws.Range("B28").Value = "username"
$ws.Range("D28").Value = "java.lang.String"
$ws.Range("B29").Value = "enabled"
$ws.Range("D29").Value = "java.lang.Boolean"
$ws.Range("B30").Value = "id"
$ws.Range("D30").Value = "java.lang.Long"
